$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update phone numbers in B2/B3 to the new number
$ws.Range("B2").Value = 5585985010594
$ws.Range("B3").Value = 5585985010594

# 2. Remove the stray E22 cell/style (far below the used data)
$ws.Range("E22").Clear()

# 3. Delete column D entirely: this removes the "LINK DO WHATSAPP" header,
#    the HYPERLINK() formulas, and all the leftover styled-but-empty D cells
#    (D4:D10), while leaving row 8 (which still holds C8) intact.
$ws.Range("D1:D10").EntireColumn.Delete()

# 3b. The old column D carried a wider "bestFit" column width; now that its
#     content lives in column C, make column C use that width (closest
#     achievable value through the ColumnWidth property).
$ws.Columns("C").ColumnWidth = 21

# 4. Re-apply the AutoFilter so it only spans A1:C1 (was A1:D1)
$ws.AutoFilterMode = $false
$ws.Range("A1:C1").AutoFilter()

# 5. Keep the workbook-level _FilterDatabase defined name in sync with the
#    new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Planilha1!`$A`$1:`$C`$1"
    }
}

# 6. The hyperlink formatting/style is no longer used anywhere, so drop the
#    now-unused "Hiperlink" cell style (mirrors Excel's own cleanup when the
#    hyperlink column is removed)
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Hiperlink") {
        $s.Delete()
    }
}

# 7. Match the final selection left by the edit
$ws.Range("B3").Select()

Write-Host "done"
